$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MaterialTable")

# Remove the VertexShader, PixelShader and GeometryShader columns (B:D) -
# the sheet now only tracks the material/shading properties.
$ws.Range("B:D").Delete()

$ws.Range("H7").Select()
